$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.125.48'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '2.374.97'
$ws.Range('E3').Value = '  +3.26%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.506'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.21'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.98%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.33'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.07%  '
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').Value = '2.746.13'
$ws.Range('E15').Value = '  +3.35%  '
$ws.Range('D16').Value = '2.376.74'
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.808'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.53%  '
$ws.Range('D18').Value = '43.092.28'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.15'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.38'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.27%  '
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.78'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E28').Value = '  +0.33%  '
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.33'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.69%  '
$ws.Range('E33').Value = '  +8.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('E35').Value = '  +4.04%  '
$ws.Range('E36').Value = '  +5.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.32'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.30'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('E39').Value = '  +4.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.14%  '
$ws.Range('E41').Value = '  -0.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -36.56%  '
$ws.Range('D43').Value = '1.958.27'
$ws.Range('E43').Value = '  +0.50%  '
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('E47').Value = '  -10.91%  '
$ws.Range('D48').Value = '2.609.78'
$ws.Range('E48').Value = '  +3.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.60'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('E50').Value = '  +1.47%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.91'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.78%  '
